$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 109.3723396666667
$ws.Range("H2").Value = 328.117019
$ws.Range("I2").Value = 0.3006244632995563
$ws.Range("J2").Value = 0.3006244632995563
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.378475333333333
$ws.Range("N2").Value = 4.135426
$ws.Range("O2").Value = 0.05609715574531157
$ws.Range("P2").Value = 0.05609715574531156
$ws.Range("Q2").Value = 150.7670723794549
$ws.Range("R2").Value = 1356.903651415094
$ws.Range("S2").Value = 0.01686417733856591
$ws.Range("T2").Value = 0.01686417733856591

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 109.3723396666667
$ws.Range("H3").Value = 328.117019
$ws.Range("I3").Value = 0.3006244632995563
$ws.Range("J3").Value = 0.3006244632995563
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.16176133333333
$ws.Range("N3").Value = 39.485284
$ws.Range("O3").Value = 0.5356188518899525
$ws.Range("P3").Value = 0.5356188518899525
$ws.Range("Q3").Value = 1439.532631160933
$ws.Range("R3").Value = 12955.7936804484
$ws.Range("S3").Value = 0.1610201298825415
$ws.Range("T3").Value = 0.1610201298825415

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 109.3723396666667
$ws.Range("H4").Value = 328.117019
$ws.Range("I4").Value = 0.3006244632995563
$ws.Range("J4").Value = 0.3006244632995563
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.829094
$ws.Range("N4").Value = 2.487282
$ws.Range("O4").Value = 0.03374004171190829
$ws.Range("P4").Value = 0.03374004171190828
$ws.Range("Q4").Value = 90.67995058359534
$ws.Range("R4").Value = 816.1195552523581
$ws.Range("S4").Value = 0.01014308193134707
$ws.Range("T4").Value = 0.01014308193134707

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 109.3723396666667
$ws.Range("H5").Value = 328.117019
$ws.Range("I5").Value = 0.3006244632995563
$ws.Range("J5").Value = 0.3006244632995563
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.203668
$ws.Range("N5").Value = 27.611004
$ws.Range("O5").Value = 0.3745439506528278
$ws.Range("P5").Value = 0.3745439506528276
$ws.Range("Q5").Value = 1006.626702675231
$ws.Range("R5").Value = 9059.640324077078
$ws.Range("S5").Value = 0.1125970741471018
$ws.Range("T5").Value = 0.1125970741471018

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 170.751104
$ws.Range("H6").Value = 512.2533120000001
$ws.Range("I6").Value = 0.4693321835689973
$ws.Range("J6").Value = 0.4693321835689973
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.378475333333333
$ws.Range("N6").Value = 4.135426
$ws.Range("O6").Value = 0.05609715574531157
$ws.Range("P6").Value = 0.05609715574531156
$ws.Range("Q6").Value = 235.3761850034347
$ws.Range("R6").Value = 2118.385665030912
$ws.Range("S6").Value = 0.0263282005979572
$ws.Range("T6").Value = 0.0263282005979572

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 170.751104
$ws.Range("H7").Value = 512.2533120000001
$ws.Range("I7").Value = 0.4693321835689973
$ws.Range("J7").Value = 0.4693321835689973
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.16176133333333
$ws.Range("N7").Value = 39.485284
$ws.Range("O7").Value = 0.5356188518899525
$ws.Range("P7").Value = 0.5356188518899525
$ws.Range("Q7").Value = 2247.385278251179
$ws.Range("R7").Value = 20226.46750426061
$ws.Range("S7").Value = 0.2513831653182308
$ws.Range("T7").Value = 0.2513831653182308

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 170.751104
$ws.Range("H8").Value = 512.2533120000001
$ws.Range("I8").Value = 0.4693321835689973
$ws.Range("J8").Value = 0.4693321835689973
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.829094
$ws.Range("N8").Value = 2.487282
$ws.Range("O8").Value = 0.03374004171190829
$ws.Range("P8").Value = 0.03374004171190828
$ws.Range("Q8").Value = 141.568715819776
$ws.Range("R8").Value = 1274.118442377984
$ws.Range("S8").Value = 0.01583528745035897
$ws.Range("T8").Value = 0.01583528745035896

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 170.751104
$ws.Range("H9").Value = 512.2533120000001
$ws.Range("I9").Value = 0.4693321835689973
$ws.Range("J9").Value = 0.4693321835689973
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.203668
$ws.Range("N9").Value = 27.611004
$ws.Range("O9").Value = 0.3745439506528278
$ws.Range("P9").Value = 0.3745439506528276
$ws.Range("Q9").Value = 1571.536471849472
$ws.Range("R9").Value = 14143.82824664525
$ws.Range("S9").Value = 0.1757855302024504
$ws.Range("T9").Value = 0.1757855302024504

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 68.96861966666667
$ws.Range("H10").Value = 206.905859
$ws.Range("I10").Value = 0.1895694499632422
$ws.Range("J10").Value = 0.1895694499632422
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.378475333333333
$ws.Range("N10").Value = 4.135426
$ws.Range("O10").Value = 0.05609715574531157
$ws.Range("P10").Value = 0.05609715574531156
$ws.Range("Q10").Value = 95.07154098454822
$ws.Range("R10").Value = 855.643868860934
$ws.Range("S10").Value = 0.01063430695914104
$ws.Range("T10").Value = 0.01063430695914104

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 68.96861966666667
$ws.Range("H11").Value = 206.905859
$ws.Range("I11").Value = 0.1895694499632422
$ws.Range("J11").Value = 0.1895694499632422
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 13.16176133333333
$ws.Range("N11").Value = 39.485284
$ws.Range("O11").Value = 0.5356188518899525
$ws.Range("P11").Value = 0.5356188518899525
$ws.Range("Q11").Value = 907.7485115421063
$ws.Range("R11").Value = 8169.736603878957
$ws.Range("S11").Value = 0.1015369711427216
$ws.Range("T11").Value = 0.1015369711427216

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 68.96861966666667
$ws.Range("H12").Value = 206.905859
$ws.Range("I12").Value = 0.1895694499632422
$ws.Range("J12").Value = 0.1895694499632422
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.829094
$ws.Range("N12").Value = 2.487282
$ws.Range("O12").Value = 0.03374004171190829
$ws.Range("P12").Value = 0.03374004171190828
$ws.Range("Q12").Value = 57.18146875391533
$ws.Range("R12").Value = 514.633218785238
$ws.Range("S12").Value = 0.006396081149063301
$ws.Range("T12").Value = 0.0063960811490633

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 68.96861966666667
$ws.Range("H13").Value = 206.905859
$ws.Range("I13").Value = 0.1895694499632422
$ws.Range("J13").Value = 0.1895694499632422
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.203668
$ws.Range("N13").Value = 27.611004
$ws.Range("O13").Value = 0.3745439506528278
$ws.Range("P13").Value = 0.3745439506528276
$ws.Range("Q13").Value = 634.7642778302708
$ws.Range("R13").Value = 5712.878500472437
$ws.Range("S13").Value = 0.07100209071231627
$ws.Range("T13").Value = 0.07100209071231625

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.72510066666667
$ws.Range("H14").Value = 44.175302
$ws.Range("I14").Value = 0.0404739031682042
$ws.Range("J14").Value = 0.04047390316820419
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.378475333333333
$ws.Range("N14").Value = 4.135426
$ws.Range("O14").Value = 0.05609715574531157
$ws.Range("P14").Value = 0.05609715574531156
$ws.Range("Q14").Value = 20.29818804985022
$ws.Range("R14").Value = 182.683692448652
$ws.Range("S14").Value = 0.00227047084964741
$ws.Range("T14").Value = 0.002270470849647409

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.72510066666667
$ws.Range("H15").Value = 44.175302
$ws.Range("I15").Value = 0.0404739031682042
$ws.Range("J15").Value = 0.04047390316820419
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 13.16176133333333
$ws.Range("N15").Value = 39.485284
$ws.Range("O15").Value = 0.5356188518899525
$ws.Range("P15").Value = 0.5356188518899525
$ws.Range("Q15").Value = 193.8082605839742
$ws.Range("R15").Value = 1744.274345255768
$ws.Range("S15").Value = 0.02167858554645865
$ws.Range("T15").Value = 0.02167858554645864

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.72510066666667
$ws.Range("H16").Value = 44.175302
$ws.Range("I16").Value = 0.0404739031682042
$ws.Range("J16").Value = 0.04047390316820419
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.829094
$ws.Range("N16").Value = 2.487282
$ws.Range("O16").Value = 0.03374004171190829
$ws.Range("P16").Value = 0.03374004171190828
$ws.Range("Q16").Value = 12.20849261212933
$ws.Range("R16").Value = 109.876433509164
$ws.Range("S16").Value = 0.001365591181138947
$ws.Range("T16").Value = 0.001365591181138946

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.72510066666667
$ws.Range("H17").Value = 44.175302
$ws.Range("I17").Value = 0.0404739031682042
$ws.Range("J17").Value = 0.04047390316820419
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.203668
$ws.Range("N17").Value = 27.611004
$ws.Range("O17").Value = 0.3745439506528278
$ws.Range("P17").Value = 0.3745439506528276
$ws.Range("Q17").Value = 135.5249378025787
$ws.Range("R17").Value = 1219.724440223208
$ws.Range("S17").Value = 0.0151592555909592
$ws.Range("T17").Value = 0.01515925559095919

